# CEDS_Code_Change_Tracker.xlsx update
# Rewrites the tail of the change-tracker table (rows 99-102 -> 99-103),
# replacing the "Review"/"In Progress" entry with the completed EDGAR work,
# and appending two brand-new completed entries describing the IO_functions.R
# .zip support and the E.UNFCCC_emissions.R rename.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reference cells already carrying the exact styles we need to reuse ---
# s="9"  -> general center/center wrap text (used throughout column A-D, H)
# s="7"  -> date format (numFmtId 14), center/center wrap text (columns E-G)
$styleRefGeneral = $ws.Range("A98")
$styleRefDate = $ws.Range("E96")

function Set-GeneralStyle($range) {
    $styleRefGeneral.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
}

function Set-DateStyle($range) {
    $styleRefDate.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
}

# --- Row 99 ---------------------------------------------------------------
$ws.Range("A99").Value2 = "Created C1.2.add_NC_emissions_EDGAR.R to process and add EDGAR default emissions data for various emissions species to their respective default emissions databases, added call to new script into C1.2.add_NC_emissions.R."
$ws.Range("B99").Value2 = 95
$ws.Range("C99").Value2 = "Jon Seibert"
$ws.Range("D99").Value2 = "Committed"

Set-DateStyle($ws.Range("E99"))
$ws.Range("E99").Value2 = 42377
Set-DateStyle($ws.Range("F99"))
$ws.Range("F99").Value2 = 42377
Set-DateStyle($ws.Range("G99"))
$ws.Range("G99").Value2 = 42377

$ws.Range("H99").Value2 = "5edb7be"
$ws.Rows.Item(99).RowHeight = 64.5

# --- Row 100 ----------------------------------------------------------------
$ws.Range("A100").Value2 = "Created NC_EDGAR_sector_mapping.csv and Master_EDGAR_sector_mapping.csv for use in module C and checking presence of all EDGAR sectors in final emissions database."
$ws.Range("B100").Value2 = 96
$ws.Range("C100").Value2 = "Jon Seibert"
$ws.Range("D100").Value2 = "Committed"
$ws.Range("E100").Value2 = 42380
$ws.Range("F100").Value2 = 42016
$ws.Range("G100").Value2 = 42017
$ws.Range("H100").Value2 = "5edb7be"
$ws.Rows.Item(100).RowHeight = 62.25

# --- Row 101 ------------------------------------------------------------
$ws.Range("A101").Value2 = "Added EDGARcheck function to analysis_functions.R"
$ws.Range("B101").Value2 = 97
$ws.Range("C101").Value2 = "Jon Seibert"
$ws.Range("D101").Value2 = "Committed"

Set-DateStyle($ws.Range("E101"))
$ws.Range("E101").Value2 = 42380
Set-DateStyle($ws.Range("F101"))
$ws.Range("F101").Value2 = 42016
Set-DateStyle($ws.Range("G101"))
$ws.Range("G101").Value2 = 42017

$ws.Range("H101").Value2 = "5edb7be"
$ws.Rows.Item(101).RowHeight = 32.25

# --- Row 102 & 103 (brand-new content) ------------------------------------
# NOTE: shared-string table entries are created in first-write order, so the
# two new descriptions (column A) must both be written *before* either "H"
# revision-number cell, to land at uniqueCount indices ...174, 175, 176 in
# the same order the target workbook uses (description, description, hash).
$ws.Range("A102").Value2 = "Upgraded IO_functions.R readData function to include ability to read one, all, or a select list of .csv files from within a .zip file. Added listZippedFiles function."
$ws.Range("A103").Value2 = "Renamed E.UNFCCC_SO2_emissions.R to E.UNFCCC_emissions.R, updated to use new readData .zip features to read all data from within large .zip files, added dummy output for species without present input data."

$ws.Range("B102").Value2 = 98
$ws.Range("C102").Value2 = "Jon Seibert"
$ws.Range("D102").Value2 = "Committed"
Set-DateStyle($ws.Range("E102"))
$ws.Range("E102").Value2 = 42384
Set-DateStyle($ws.Range("F102"))
$ws.Range("F102").Value2 = 42387
Set-DateStyle($ws.Range("G102"))
$ws.Range("G102").Value2 = 42387
$ws.Range("H102").Value2 = "ce6f6a3"
$ws.Rows.Item(102).RowHeight = 54.75

$ws.Range("B103").Value2 = 99
$ws.Range("C103").Value2 = "Jon Seibert"
$ws.Range("D103").Value2 = "Committed"
Set-DateStyle($ws.Range("E103"))
$ws.Range("E103").Value2 = 42384
Set-DateStyle($ws.Range("F103"))
$ws.Range("F103").Value2 = 42387
Set-DateStyle($ws.Range("G103"))
$ws.Range("G103").Value2 = 42387
$ws.Range("H103").Value2 = "ce6f6a3"
$ws.Rows.Item(103).RowHeight = 73.5

# --- Row 104 (new, blank spacer row) --------------------------------------
$ws.Rows.Item(104).RowHeight = 36.75

# --- Selection / active cell ----------------------------------------------
$ws.Range("C105").Select()
